$wb = $excel.ActiveWorkbook

# ALC row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3471.1702
$ws.Range("I76").Value = 3026.442
$ws.Range("K76").Value = 3026.442
$ws.Range("M76").Value = -2711.442

# ALC row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3471.1702
$ws.Range("I79").Value = 3026.442
$ws.Range("K79").Value = 3026.442
$ws.Range("M79").Value = -1934.442

# ARM row 2: Ain't Got No Ingots | Bronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1247.8462
$ws.Range("J2").Value = 1560
$ws.Range("L2").Value = 1560
$ws.Range("N2").Value = -1786

# ARM row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2788.0557
$ws.Range("I61").Value = 1636.5625
$ws.Range("J61").Value = 12000
$ws.Range("K61").Value = 1636.5625
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = -1424.5625
$ws.Range("N61").Value = -12424

# ARM row 74: As the Bolt Flies | Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 49809.12
$ws.Range("I74").Value = 56504.86
$ws.Range("J74").Value = 1599.8
$ws.Range("K74").Value = 56504.86
$ws.Range("L74").Value = 1599.8
$ws.Range("M74").Value = -55630.86
$ws.Range("N74").Value = -3347.8

# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 49809.12
$ws.Range("I77").Value = 56504.86
$ws.Range("J77").Value = 1599.8
$ws.Range("K77").Value = 282524.3
$ws.Range("L77").Value = 7999
$ws.Range("M77").Value = -278156.3
$ws.Range("N77").Value = -16735

# ARM row 88: The Mast Chance | Adamantite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 15611.866
$ws.Range("I88").Value = 1980
$ws.Range("J88").Value = 22427.8
$ws.Range("K88").Value = 1980
$ws.Range("L88").Value = 22427.8
$ws.Range("M88").Value = -1574
$ws.Range("N88").Value = -23239.8

# ARM row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 15611.866
$ws.Range("I91").Value = 1980
$ws.Range("J91").Value = 22427.8
$ws.Range("K91").Value = 1980
$ws.Range("L91").Value = 22427.8
$ws.Range("M91").Value = -576
$ws.Range("N91").Value = -25235.8

# ARM row 116: No Scope | Titanbronze Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1247.8462
$ws.Range("J116").Value = 1560
$ws.Range("L116").Value = 1560
$ws.Range("N116").Value = -6148

# ARM row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2788.0557
$ws.Range("I136").Value = 1636.5625
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 4909.6875
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -2359.6875
$ws.Range("N136").Value = -41100

# BSM row 3: Hells Bells | Bronze Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1247.8462
$ws.Range("J3").Value = 1560
$ws.Range("L3").Value = 1560
$ws.Range("N3").Value = -1788

# CRP row 31: Wall Not Found | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18224.42
$ws.Range("I31").Value = 27920.352
$ws.Range("J31").Value = 3874.44
$ws.Range("K31").Value = 27920.352
$ws.Range("L31").Value = 3874.44
$ws.Range("M31").Value = -27625.352
$ws.Range("N31").Value = -4464.440000000001

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 18224.42
$ws.Range("I34").Value = 27920.352
$ws.Range("J34").Value = 3874.44
$ws.Range("K34").Value = 27920.352
$ws.Range("L34").Value = 3874.44
$ws.Range("M34").Value = -27718.352
$ws.Range("N34").Value = -4278.440000000001

# CRP row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2677.2058
$ws.Range("I58").Value = 977.13635
$ws.Range("J58").Value = 5794
$ws.Range("K58").Value = 977.13635
$ws.Range("L58").Value = 5794
$ws.Range("M58").Value = -774.13635
$ws.Range("N58").Value = -6200

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1753.3
$ws.Range("I132").Value = 1233.381
$ws.Range("J132").Value = 2966.4443
$ws.Range("K132").Value = 3700.143
$ws.Range("L132").Value = 8899.332900000001
$ws.Range("M132").Value = -1170.143
$ws.Range("N132").Value = -13959.3329

# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11365017
$ws.Range("I134").Value = 1147.6857
$ws.Range("J134").Value = 55557840
$ws.Range("K134").Value = 3443.0571
$ws.Range("L134").Value = 166673520
$ws.Range("M134").Value = -908.0571
$ws.Range("N134").Value = -166678590

# CRP row 136: Turali Quality | Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2677.2058
$ws.Range("I136").Value = 977.13635
$ws.Range("J136").Value = 5794
$ws.Range("K136").Value = 2931.40905
$ws.Range("L136").Value = 17382
$ws.Range("M136").Value = -381.4090500000002
$ws.Range("N136").Value = -22482

# CUL row 22: A Total Nut Job | Walnut Bread
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1146.0769
$ws.Range("I22").Value = 1016.5
$ws.Range("J22").Value = 1257.1428
$ws.Range("K22").Value = 3049.5
$ws.Range("L22").Value = 3771.4284
$ws.Range("M22").Value = -2880.5
$ws.Range("N22").Value = -4109.428400000001

# CUL row 27: Brain Food | Walnut Bread
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 1146.0769
$ws.Range("I27").Value = 1016.5
$ws.Range("J27").Value = 1257.1428
$ws.Range("K27").Value = 3049.5
$ws.Range("L27").Value = 3771.4284
$ws.Range("M27").Value = -2947.5
$ws.Range("N27").Value = -3975.4284

# CUL row 113: Can't Eat Just One | Night Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 595.9286
$ws.Range("I113").Value = 437.75
$ws.Range("J113").Value = 659.2
$ws.Range("K113").Value = 1313.25
$ws.Range("L113").Value = 1977.6
$ws.Range("M113").Value = 856.75
$ws.Range("N113").Value = -6317.6

# GSM row 15: The Tusk at Hand | Fang Earrings
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19227.285
$ws.Range("J15").Value = 19227.285
$ws.Range("L15").Value = 19227.285
$ws.Range("N15").Value = -19803.285

# GSM row 81: The Grander Temple | Dragon Fang Earrings
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 19227.285
$ws.Range("J81").Value = 19227.285
$ws.Range("L81").Value = 19227.285
$ws.Range("N81").Value = -21223.285

# GSM row 84: Man with a Dragon Earring (L) | Dragon Fang Earrings
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 19227.285
$ws.Range("J84").Value = 19227.285
$ws.Range("L84").Value = 57681.855
$ws.Range("N84").Value = -67665.855

# GSM row 127: Sage with the Golden Earrings | Phrygian Ear Cuffs of Healing
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 38771
$ws.Range("J127").Value = 38771
$ws.Range("L127").Value = 38771
$ws.Range("N127").Value = -48691

# LTW row 40: Best Served Toad | Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 46352.39
$ws.Range("I40").Value = 1865.8334
$ws.Range("J40").Value = 94883.17999999999
$ws.Range("K40").Value = 1865.8334
$ws.Range("L40").Value = 94883.17999999999
$ws.Range("M40").Value = -1729.8334
$ws.Range("N40").Value = -95155.17999999999

# LTW row 46: Supply Side Logic | Boar Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1933.6666
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 2450.5
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 2450.5
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -1926.25

# WVR row 30: The Telltale Tress | Cotton Coif of Gathering
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 70010
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 70010
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 70010
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -70224

# WVR row 74: Clothing the Naked Truth | Ramie Robe of Casting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8699.857
$ws.Range("J74").Value = 8699.857
$ws.Range("L74").Value = 8699.857
$ws.Range("N74").Value = -10571.857

# WVR row 77: When in Robes (L) | Ramie Robe of Casting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 8699.857
$ws.Range("J77").Value = 8699.857
$ws.Range("L77").Value = 26099.571
$ws.Range("N77").Value = -35459.571
